$d = $word.ActiveDocument

# 1) Fix the paragraph-mark/run color of the "¿Cuántos empleados..." paragraph
#    (item 22) from theme-based black (000000/text1) to plain 202124.
$p22 = $d.Paragraphs.Item(22)
$p22.Range.Font.Color = 2367776

# 2) Remove the block of "process/system" interview questions (items 23-40),
#    i.e. everything from "A quien o a quienes afecta directamente." through
#    "Que servicios espera que provea.", collapsing it down to a single blank
#    paragraph (which is what remains of paragraph 41's predecessor, the
#    trailing blank paragraph that absorbs the merge).
$pStart = $d.Paragraphs.Item(23)
$pEnd = $d.Paragraphs.Item(40)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$rng.Delete()

# 3) The remaining blank paragraph needs a left indent of 720 twips (36 pt)
#    added to its paragraph formatting.
$pBlank = $d.Paragraphs.Item(23)
$pBlank.Range.ParagraphFormat.LeftIndent = 36
